$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    3  = @{ BR = 70.2045;             BS = 14.75;              BT = 1.795500000000004 }
    4  = @{ BR = 69.81399999999999;   BS = 13.5;               BT = 2.186000000000007 }
    10 = @{ BR = 69.92325;            BS = 14.375;             BT = 2.076750000000004 }
    13 = @{ BR = 69.65100000000001;   BS = 8.333333333333334;  BT = 2.34899999999999 }
    17 = @{ BR = 71.15966666666667;   BS = 38.5;               BT = 0.8403333333333336 }
    19 = @{ BR = 70.40066666666667;   BS = 50.5;               BT = 1.599333333333334 }
    23 = @{ BR = 70.52124999999999;   BS = 20.125;             BT = 1.478750000000005 }
    30 = @{ BR = 70.815;              BS = 27.125;             BT = 1.185000000000002 }
    32 = @{ BR = 71.72525;            BS = 66.5;               BT = 0.2747499999999974 }
    34 = @{ BR = 69.46899999999999;   BS = 33.75;              BT = 2.531000000000006 }
    35 = @{ BR = 70.62625;            BS = 23.5;               BT = 1.373750000000001 }
    36 = @{ BR = 71.69199999999999;   BS = 72.625;             BT = 0.3080000000000069 }
    39 = @{ BR = 69.91733333333333;   BS = 49.83333333333334;  BT = 2.082666666666668 }
    40 = @{ BR = 70.71225;            BS = 53;                 BT = 1.287750000000003 }
    45 = @{ BR = 70.37725;            BS = 54.375;             BT = 1.622749999999996 }
    48 = @{ BR = 70.53175;            BS = 29.25;              BT = 1.468249999999998 }
    59 = @{ BR = 71.69924999999999;   BS = 74.375;             BT = 0.3007500000000078 }
    65 = @{ BR = 71.66799999999999;   BS = 75;                 BT = 0.3320000000000078 }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Range("BR$row").Value = $vals.BR
    $ws.Range("BS$row").Value = $vals.BS
    $ws.Range("BT$row").Value = $vals.BT
}
